$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-08 Friday" "2025-08-09 Saturday"

Replace-Text "590×5=2950" "407×4=1628"
Replace-Text "648×7=4536" "921×7=6447"
Replace-Text "637×6=3822" "997×7=6979"
Replace-Text "311×8=2488" "253×7=1771"
Replace-Text "339×9=3051" "453×2=906"

Replace-Text "693×6=4158" "219×9=1971"
Replace-Text "424×7=2968" "723×5=3615"
Replace-Text "260×7=1820" "668×2=1336"
Replace-Text "416×4=1664" "758×8=6064"
Replace-Text "370×9=3330" "205×3=615"

Replace-Text "138×7=966" "652×6=3912"
Replace-Text "239×2=478" "255×9=2295"
Replace-Text "150×6=900" "759×5=3795"
Replace-Text "566×7=3962" "782×7=5474"
Replace-Text "767×6=4602" "506×2=1012"

Replace-Text "276×9=2484" "556×4=2224"
Replace-Text "854×3=2562" "946×8=7568"
Replace-Text "733×2=1466" "665×4=2660"
Replace-Text "921×5=4605" "773×5=3865"
Replace-Text "516×8=4128" "823×6=4938"

Replace-Text "259×7=1813" "509×6=3054"
Replace-Text "378×4=1512" "606×5=3030"
Replace-Text "638×6=3828" "893×4=3572"
Replace-Text "293×3=879" "744×4=2976"
Replace-Text "622×6=3732" "154×4=616"
